$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.201.58"
$ws.Range("E2").Value = "  +0.32%  "

$ws.Range("D3").Value = "1.649.77"
$ws.Range("E3").Value = "  +0.14%  "

$ws.Range("E4").Value = "  +0.07%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "218.74"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.54%  "

$ws.Range("E6").Value = "  +2.23%  "

$ws.Range("E7").Value = "  +0.11%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.256"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.05%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.0628"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.35%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "20.21"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +3.04%  "

$ws.Range("E11").Value = "  -0.08%  "

$ws.Range("D12").Value = "1.882.71"
$ws.Range("E12").Value = "  +0.35%  "

$ws.Range("D13").Value = "1.645.34"
$ws.Range("E13").Value = "  -0.20%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "4.14"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.33%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.538"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.09%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "67.92"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.58%  "

$ws.Range("D17").Value = "27.191.36"
$ws.Range("E17").Value = "  +0.55%  "

$ws.Range("D18").Value = "0.0₃0737"
$ws.Range("E18").Value = "  +0.30%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "220.51"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("E20").Value = "  -0.08%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.73"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.16%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.44"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.48%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "2.48"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.63%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "9.26"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.00%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "148.32"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.44%  "

$ws.Range("E26").Value = "  +0.00%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "7.39"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.61%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.119"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.66%  "

$ws.Range("E29").Value = "  -0.11%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.0507"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.35%  "

$ws.Range("E31").Value = "  -0.32%  "

$ws.Range("E32").Value = "  -0.70%  "

$ws.Range("E33").Value = "  +0.58%  "

$ws.Range("E34").Value = "  +0.47%  "

$ws.Range("D35").Value = "1.272.25"
$ws.Range("E35").Value = "  +0.53%  "

$ws.Range("E36").Value = "  +1.30%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.0177"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.22%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.541"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.27%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.844"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +2.01%  "

$ws.Range("E40").Value = "  +0.04%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.810"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.39%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "5.41"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.61%  "

$ws.Range("D43").Value = "1.792.90"
$ws.Range("E43").Value = "  +0.36%  "

$ws.Range("E44").Value = "  +5.22%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "63.21"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.99%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "92.41"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.28%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.59"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.84%  "

$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("E48").Value = "  +14.81%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0514"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.27%  "

$ws.Range("E50").Value = "  +1.05%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0976"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.38%  "
